$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -19.05088809060817
$ws.Range("K2").Value = -7.333333333333343

$ws.Range("B3").Value = -1.777663414725339
$ws.Range("K3").Value = -7.333333333333343

$ws.Range("B4").Value = 0.1522293178118161
$ws.Range("K4").Value = 3.666666666666657

$ws.Range("B5").Value = 2.221169754193852
$ws.Range("K5").Value = 7.333333333333343

$ws.Range("B6").Value = 4.289544203869809
$ws.Range("I6").Value = -36.6863929236917
$ws.Range("K6").Value = 7.333333333333343

$ws.Range("B7").Value = 15.26087022847696
$ws.Range("I7").Value = -36.6863929236917
$ws.Range("K7").Value = -3.666666666666657

$ws.Range("B8").Value = 17.93987995906764
$ws.Range("K8").Value = -7.333333333333343

$ws.Range("B9").Value = 10.13649656216057
$ws.Range("K9").Value = -11

$ws.Range("B10").Value = -19.89296983070199
$ws.Range("K10").Value = -7.333333333333343

$ws.Range("B11").Value = -52.44123335746599

$ws.Range("B12").Value = -64.00719282167165

$ws.Range("B13").Value = -67.25603530169883
$ws.Range("I13").Value = -37.23736614058328

$ws.Range("B14").Value = -63.56512750605589
$ws.Range("I14").Value = -36.741710257504

$ws.Range("B15").Value = -59.59222878126275

$ws.Range("B16").Value = -58.03568624059153
$ws.Range("I16").Value = -36.72505219675372

$ws.Range("B17").Value = -56.60720464534006

$ws.Range("B18").Value = -49.53998385111981
$ws.Range("I18").Value = -37.56895584042388

$ws.Range("B19").Value = -46.36551164545335
$ws.Range("I19").Value = -37.98697887434616

$ws.Range("B20").Value = -53.89163420216983
$ws.Range("I20").Value = -37.90494578272683

$ws.Range("B21").Value = -72.51460329834299

$ws.Range("B22").Value = -84.1261813121655

$ws.Range("B23").Value = -92.24986939425389

$ws.Range("B24").Value = -94.89035873264947

$ws.Range("B25").Value = -76.27470972491619
